# "combine bees and butterflies into butterflies script"
# For a set of plantfamily rows, the croppedarea (column F) value changes
# from the text "NA" to the number 1, and the cell is highlighted with a
# green fill (theme accent6) to flag the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(10, 23, 26, 27, 30, 31, 34, 40, 43, 48, 51, 56)
foreach ($r in $rows) {
    $cell = $ws.Range("F$r")
    $cell.Value = 1
    $cell.Interior.ThemeColor = 10
}

# Update the view state: scroll down and select F57 (best effort - the
# scroll position recorded in the pane is a cosmetic view setting).
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("F57").Select()

Write-Host "Updated croppedarea cells for 12 plant families"
